$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: DRA5 row gets new Jira id / Description text (wrapped style, like C6/C7) ---
$ws.Range("C6").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("C6").Copy()
$ws.Range("C10").PasteSpecial(-4122)

$ws.Range("B10").Value = "OPQA-4205||OPQA-4207||OPQA-4208||OPQA-4210||OPQA-4211"
$ws.Range("C10").Value = "Verify that the profile fly-out should contain link to terms of use||Verify that profile fly-out should contain link to privacy statement||Verify that the profile fly-out should contain link to app-specific feedback page (http://thomson-reuters1260211048.drug-research-advisor-target-druggability-feedback.sgizmo.com/s3)||Verify that the profile fly-out should contain link to app-specific help page||Verify that the alternative profile fly-out should contain link to sign out of the platform. User returns to DRA sign-in page."

$ws.Rows.Item(10).RowHeight = 90

# --- Row 11: brand new test case DRA6 ---
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C6").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("E9").Copy()
$ws.Range("E11").PasteSpecial(-4122)

$ws.Range("A11").Value = "DRA6"
$ws.Range("B11").Value = "OPQA-4197||OPQA-4199||OPQA-4215||OPQA-4216"
$ws.Range("C11").Value = "Verify that profile fly-out will display profile meta-data||Verify that profile fly-out provides access to the profile modal.||Verify that the profile fly-out should display the following user profile details, if available: a)First name b)Last Name c)Title d)Institution e)Country f)Photo||Verify that by clicking on any of the following fields (when present), will provide access to the profile modal. 1.Name 2.Institution 3. Country 4 .Title 5.Photo"
$ws.Range("D11").Value = "Y"

$ws.Rows.Item(11).RowHeight = 75

# --- Sheet view: scrolled/selected further down after adding the new row ---
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("C16:C17").Select()
